$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.530.77"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "3.679.73"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "643.35"
$ws.Range("E5").Value = "  -5.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.81"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.498"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.11"
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "4.295.50"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.74"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "3.688.48"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "69.479.13"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.01"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.00"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.88"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.649"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.36"
$ws.Range("E23").Value = "  -0.87%  "
$ws.Range("D24").Value = "3.825.19"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.93"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.04"
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.63"
$ws.Range("E29").Value = "  -2.58%  "
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.95"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.47"
$ws.Range("E34").Value = "  -1.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.165"
$ws.Range("E35").Value = "  +3.61%  "
$ws.Range("D36").Value = "3.672.83"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.48"
$ws.Range("E37").Value = "  +1.89%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.90"
$ws.Range("E39").Value = "  -6.80%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "177.35"
$ws.Range("E41").Value = "  +4.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.21"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0901"
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.88"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.56"
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("E48").Value = "  -4.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.08"
$ws.Range("E49").Value = "  -2.99%  "
$ws.Range("E50").Value = "  -3.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.86"
$ws.Range("E51").Value = "  +0.36%  "
